# Auto-generated Excel COM-interop script to apply the 2026-01-16 16:14:21 scrape update
$wb = $excel.ActiveWorkbook

$updated = 'Última actualización: 16:14:21'

# --- Sheet "LP1912" ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = $updated
$ws1.Range("A3").Value = 'Total filas: 373'
$data1 = @(
  @('00:04:05','01:12','215_ALUAR',68,'LP1912'),
  @('00:04:05','01:58','14_ABASTO',114,'LP1912'),
  @('01:07:17','02:58','215_ALUAR',111,'LP1912'),
  @('01:49:17','03:03','215_ALUAR',74,'LP1912'),
  @('01:49:17','03:48','14_ABASTO',119,'LP1912'),
  @('02:42:39','03:54','14_ABASTO',72,'LP1912'),
  @('02:13:28','03:55','14_ABASTO',102,'LP1912'),
  @('02:13:28','04:01','81_EL PELIGRO',108,'LP1912'),
  @('03:23:38','04:45','215A_EL PATO',82,'LP1912'),
  @('02:56:05','04:46','215A_EL PATO',110,'LP1912'),
  @('02:56:05','04:53','11_ETCHEVERRY',117,'LP1912'),
  @('04:53:50','04:54','11_ETCHEVERRY',1,'LP1912'),
  @('04:53:50','05:15','14_ABASTO',22,'LP1912'),
  @('03:23:38','05:16','17_ROMERO',113,'LP1912'),
  @('04:40:32','05:17','14_ABASTO',37,'LP1912'),
  @('05:18:23','05:20','14_ABASTO',2,'LP1912'),
  @('03:23:38','05:22','23_HERNANDEZ',119,'LP1912'),
  @('03:54:40','05:34','215B_EL PATO',100,'LP1912'),
  @('04:18:06','05:35','215B_EL PATO',77,'LP1912'),
  @('04:18:06','05:36','14_ABASTO',78,'LP1912'),
  @('03:54:40','05:39','14_ABASTO',105,'LP1912'),
  @('03:54:40','05:46','15_ABASTO',112,'LP1912'),
  @('05:47:32','05:47','15_ABASTO',0,'LP1912'),
  @('05:47:32','05:49','14_ABASTO',2,'LP1912'),
  @('04:40:32','06:04','16_SANTA ANA',84,'LP1912'),
  @('06:02:16','06:05','16_SANTA ANA',3,'LP1912'),
  @('04:18:06','06:09','16_SANTA ANA',111,'LP1912'),
  @('04:40:32','06:11','215A_EL PATO',91,'LP1912'),
  @('04:18:06','06:12','215A_EL PATO',114,'LP1912'),
  @('04:18:06','06:14','225_HARAS DEL SUR',116,'LP1912'),
  @('04:40:32','06:21','26_HERNANDEZ',101,'LP1912'),
  @('04:40:32','06:27','23_HERNANDEZ',107,'LP1912'),
  @('04:40:32','06:29','86_EST CHICA-ESC AGRARIA',109,'LP1912'),
  @('06:02:16','06:30','86_EST CHICA-ESC AGRARIA',28,'LP1912'),
  @('04:40:32','06:31','16_SANTA ANA',111,'LP1912'),
  @('04:53:50','06:44','225_C ROCA-H SUR',111,'LP1912'),
  @('04:53:50','06:46','215C_EL PATO',113,'LP1912'),
  @('05:18:23','06:58','10_OLMOS',100,'LP1912'),
  @('05:18:23','06:59','14_ABASTO',101,'LP1912'),
  @('06:02:16','07:00','14_ABASTO',58,'LP1912'),
  @('06:37:24','07:01','16_SANTA ANA',24,'LP1912'),
  @('05:47:32','07:04','23_HERNANDEZ',77,'LP1912'),
  @('05:18:23','07:05','15_ABASTO',107,'LP1912'),
  @('06:02:16','07:05','23_HERNANDEZ',63,'LP1912'),
  @('05:18:23','07:07','225_GOMEZ',109,'LP1912'),
  @('05:18:23','07:11','215A_EL PATO',113,'LP1912'),
  @('07:14:27','07:14','11_ETCHEVERRY',0,'LP1912'),
  @('05:18:23','07:15','11_ETCHEVERRY',117,'LP1912'),
  @('06:37:24','07:16','16_SANTA ANA',39,'LP1912'),
  @('06:02:16','07:16','11_ETCHEVERRY',74,'LP1912'),
  @('05:47:32','07:21','26_HERNANDEZ',94,'LP1912'),
  @('06:02:16','07:23','10_OLMOS',81,'LP1912'),
  @('05:47:32','07:27','10_OLMOS',100,'LP1912'),
  @('05:47:32','07:31','16_SANTA ANA',104,'LP1912'),
  @('05:47:32','07:31','11_ETCHEVERRY',104,'LP1912'),
  @('05:47:32','07:32','84_COLONIA URQUIZA-ESC 49',105,'LP1912'),
  @('06:02:16','07:32','11_ETCHEVERRY',90,'LP1912'),
  @('06:37:24','07:34','23_HERNANDEZ',57,'LP1912'),
  @('05:47:32','07:36','27_EL RETIRO',109,'LP1912'),
  @('07:14:27','07:37','23_HERNANDEZ',23,'LP1912'),
  @('06:02:16','07:37','27_EL RETIRO',95,'LP1912'),
  @('05:47:32','07:39','10_OLMOS',112,'LP1912'),
  @('07:44:08','07:45','14_ABASTO',1,'LP1912'),
  @('06:37:24','07:47','14_ABASTO',70,'LP1912'),
  @('06:02:16','07:48','14_ABASTO',106,'LP1912'),
  @('06:02:16','07:51','215D_EL PATO',109,'LP1912'),
  @('07:44:08','07:52','215D_EL PATO',8,'LP1912'),
  @('07:44:08','07:55','10_OLMOS',11,'LP1912'),
  @('07:14:27','07:58','16_SANTA ANA',44,'LP1912'),
  @('07:14:27','08:03','11_ETCHEVERRY',49,'LP1912'),
  @('07:44:08','08:04','11_ETCHEVERRY',20,'LP1912'),
  @('07:57:27','08:06','11_ETCHEVERRY',9,'LP1912'),
  @('07:57:27','08:10','16_SANTA ANA',13,'LP1912'),
  @('07:44:08','08:11','16_SANTA ANA',27,'LP1912'),
  @('06:37:24','08:12','15_ABASTO',95,'LP1912'),
  @('07:44:08','08:13','10_OLMOS',29,'LP1912'),
  @('06:37:24','08:21','26_HERNANDEZ',104,'LP1912'),
  @('06:37:24','08:22','16_P MOR-SANTA ANA',105,'LP1912'),
  @('06:37:24','08:23','215B_EL PATO',106,'LP1912'),
  @('07:44:08','08:23','16_P MOR-SANTA ANA',39,'LP1912'),
  @('06:37:24','08:27','84_COLONIA URQUIZA-ESC 49',110,'LP1912'),
  @('07:57:27','08:33','10_OLMOS',36,'LP1912'),
  @('06:52:38','08:42','81_EL PELIGRO',110,'LP1912'),
  @('07:14:27','08:43','14_ABASTO',89,'LP1912'),
  @('07:44:08','08:44','14_ABASTO',60,'LP1912'),
  @('07:14:27','08:54','17_ROMERO',100,'LP1912'),
  @('08:47:19','08:54','10_OLMOS',7,'LP1912'),
  @('08:16:48','08:55','10_OLMOS',39,'LP1912'),
  @('07:14:27','09:01','215A_EL PATO',107,'LP1912'),
  @('07:44:08','09:02','215A_EL PATO',78,'LP1912'),
  @('07:57:27','09:03','11_ETCHEVERRY',66,'LP1912'),
  @('08:16:48','09:04','23_HERNANDEZ',48,'LP1912'),
  @('08:47:19','09:05','23_HERNANDEZ',18,'LP1912'),
  @('07:14:27','09:07','23_HERNANDEZ',113,'LP1912'),
  @('07:44:08','09:08','23_HERNANDEZ',84,'LP1912'),
  @('07:57:27','09:09','23_HERNANDEZ',72,'LP1912'),
  @('07:14:27','09:10','16_P MOR-SANTA ANA',116,'LP1912'),
  @('07:44:08','09:11','16_P MOR-SANTA ANA',87,'LP1912'),
  @('08:33:47','09:13','10_OLMOS',40,'LP1912'),
  @('07:44:08','09:14','16_SANTA ANA',90,'LP1912'),
  @('07:57:27','09:16','27_EL RETIRO',79,'LP1912'),
  @('07:44:08','09:17','27_EL RETIRO',93,'LP1912'),
  @('07:44:08','09:21','26_HERNANDEZ',97,'LP1912'),
  @('07:57:27','09:22','16_SANTA ANA',85,'LP1912'),
  @('07:57:27','09:22','17_ROMERO',85,'LP1912'),
  @('08:47:19','09:23','16_SANTA ANA',36,'LP1912'),
  @('07:44:08','09:23','17_ROMERO',99,'LP1912'),
  @('07:57:27','09:23','11_ETCHEVERRY',86,'LP1912'),
  @('07:44:08','09:24','11_ETCHEVERRY',100,'LP1912'),
  @('08:16:48','09:29','16_SANTA ANA',73,'LP1912'),
  @('07:44:08','09:32','15_ABASTO',108,'LP1912'),
  @('07:44:08','09:33','10_OLMOS',109,'LP1912'),
  @('08:54:42','09:34','23_HERNANDEZ',40,'LP1912'),
  @('08:33:47','09:34','16_SANTA ANA',61,'LP1912'),
  @('08:47:19','09:35','16_SANTA ANA',48,'LP1912'),
  @('08:47:19','09:35','23_HERNANDEZ',48,'LP1912'),
  @('07:44:08','09:36','23_HERNANDEZ',112,'LP1912'),
  @('08:16:48','09:37','23_HERNANDEZ',81,'LP1912'),
  @('08:33:47','09:41','23_HERNANDEZ',68,'LP1912'),
  @('08:16:48','09:41','215C_EL PATO',85,'LP1912'),
  @('07:44:08','09:42','215C_EL PATO',118,'LP1912'),
  @('07:57:27','09:43','14_ABASTO',106,'LP1912'),
  @('09:25:30','09:46','16_SANTA ANA',21,'LP1912'),
  @('08:54:42','09:52','15_ABASTO',58,'LP1912'),
  @('08:54:42','09:53','10_OLMOS',59,'LP1912'),
  @('09:25:30','09:58','16_SANTA ANA',33,'LP1912'),
  @('09:25:30','10:03','11_ETCHEVERRY',38,'LP1912'),
  @('09:25:30','10:04','23_HERNANDEZ',39,'LP1912'),
  @('08:16:48','10:10','16_P MOR-SANTA ANA',114,'LP1912'),
  @('10:11:11','10:11','16_P MOR-SANTA ANA',0,'LP1912'),
  @('08:16:48','10:12','15_ABASTO',116,'LP1912'),
  @('09:25:30','10:13','10_OLMOS',48,'LP1912'),
  @('08:33:47','10:21','26_HERNANDEZ',108,'LP1912'),
  @('08:33:47','10:22','17_ROMERO',109,'LP1912'),
  @('10:11:11','10:23','16_SANTA ANA',12,'LP1912'),
  @('09:25:30','10:23','11_ETCHEVERRY',58,'LP1912'),
  @('10:11:11','10:24','11_ETCHEVERRY',13,'LP1912'),
  @('08:33:47','10:26','215A_EL PATO',113,'LP1912'),
  @('10:11:11','10:27','215A_EL PATO',16,'LP1912'),
  @('10:11:11','10:33','10_OLMOS',22,'LP1912'),
  @('10:11:11','10:36','23_HERNANDEZ',25,'LP1912'),
  @('08:54:42','10:41','17_ROMERO',107,'LP1912'),
  @('08:47:19','10:42','17_ROMERO',115,'LP1912'),
  @('08:47:19','10:43','14_ABASTO',116,'LP1912'),
  @('10:11:11','10:44','14_ABASTO',33,'LP1912'),
  @('10:11:11','10:47','16_SANTA ANA',36,'LP1912'),
  @('10:50:37','10:51','16_SANTA ANA',1,'LP1912'),
  @('10:11:11','10:52','15_ABASTO',41,'LP1912'),
  @('09:25:30','10:53','27_EL RETIRO',88,'LP1912'),
  @('10:50:37','10:56','27_EL RETIRO',6,'LP1912'),
  @('10:11:11','10:57','27_EL RETIRO',46,'LP1912'),
  @('10:50:37','11:01','215C_EL PATO',11,'LP1912'),
  @('09:25:30','11:02','215C_EL PATO',97,'LP1912'),
  @('10:50:37','11:03','11_ETCHEVERRY',13,'LP1912'),
  @('10:50:37','11:04','23_HERNANDEZ',14,'LP1912'),
  @('10:11:11','11:04','11_ETCHEVERRY',53,'LP1912'),
  @('10:11:11','11:05','23_HERNANDEZ',54,'LP1912'),
  @('09:25:30','11:06','16_P MOR-167 Y 521',101,'LP1912'),
  @('10:11:11','11:07','16_P MOR-167 Y 521',56,'LP1912'),
  @('10:11:11','11:11','10_OLMOS',60,'LP1912'),
  @('10:11:11','11:12','15_ABASTO',61,'LP1912'),
  @('11:15:53','11:15','86_EST CHICA-ESC AGRARIA',0,'LP1912'),
  @('11:15:53','11:15','16_SANTA ANA',0,'LP1912'),
  @('11:15:53','11:16','15_ABASTO',1,'LP1912'),
  @('09:25:30','11:19','86_EST CHICA-ESC AGRARIA',114,'LP1912'),
  @('10:11:11','11:20','86_EST CHICA-ESC AGRARIA',69,'LP1912'),
  @('11:15:53','11:20','26_HERNANDEZ',5,'LP1912'),
  @('09:25:30','11:21','26_HERNANDEZ',116,'LP1912'),
  @('11:15:53','11:26','225_C ROCA-H SUR',11,'LP1912'),
  @('10:11:11','11:27','225_C ROCA-H SUR',76,'LP1912'),
  @('10:11:11','11:32','81_EL PELIGRO',81,'LP1912'),
  @('10:50:37','11:34','23_HERNANDEZ',44,'LP1912'),
  @('10:50:37','11:35','11_ETCHEVERRY',45,'LP1912'),
  @('11:15:53','11:35','23_HERNANDEZ',20,'LP1912'),
  @('10:11:11','11:38','10_OLMOS',87,'LP1912'),
  @('10:50:37','11:41','17_ROMERO',51,'LP1912'),
  @('10:11:11','11:42','17_ROMERO',91,'LP1912'),
  @('11:43:19','11:43','16_SANTA ANA',0,'LP1912'),
  @('11:43:19','11:43','17_ROMERO',0,'LP1912'),
  @('10:50:37','11:43','10_OLMOS',53,'LP1912'),
  @('11:43:19','11:44','16_SANTA ANA',1,'LP1912'),
  @('10:11:11','11:51','215B_EL PATO',100,'LP1912'),
  @('11:15:53','11:58','225_GOMEZ',43,'LP1912'),
  @('11:58:46','11:58','17_ROMERO',0,'LP1912'),
  @('11:58:46','11:59','16_SANTA ANA',1,'LP1912'),
  @('10:11:11','11:59','225_GOMEZ',108,'LP1912'),
  @('11:58:46','12:00','225_GOMEZ',2,'LP1912'),
  @('10:11:11','12:02','84_COLONIA URQUIZA-ESC 49',111,'LP1912'),
  @('11:15:53','12:04','23_HERNANDEZ',49,'LP1912'),
  @('10:50:37','12:06','16_P MOR-SANTA ANA',76,'LP1912'),
  @('10:50:37','12:06','14_ABASTO',76,'LP1912'),
  @('10:50:37','12:07','10_OLMOS',77,'LP1912'),
  @('10:11:11','12:07','16_P MOR-SANTA ANA',116,'LP1912'),
  @('10:11:11','12:07','14_ABASTO',116,'LP1912'),
  @('11:43:19','12:09','23_HERNANDEZ',26,'LP1912'),
  @('11:15:53','12:12','10_OLMOS',57,'LP1912'),
  @('11:43:19','12:13','10_OLMOS',30,'LP1912'),
  @('11:15:53','12:20','14_ABASTO',65,'LP1912'),
  @('11:15:53','12:20','26_HERNANDEZ',65,'LP1912'),
  @('10:50:37','12:20','215A_EL PATO',90,'LP1912'),
  @('10:50:37','12:21','14_ABASTO',91,'LP1912'),
  @('10:50:37','12:21','26_HERNANDEZ',91,'LP1912'),
  @('10:50:37','12:22','17_ROMERO',92,'LP1912'),
  @('12:24:14','12:24','17_ROMERO',0,'LP1912'),
  @('12:24:14','12:24','16_SANTA ANA',0,'LP1912'),
  @('12:24:14','12:27','10_OLMOS',3,'LP1912'),
  @('11:58:46','12:34','23_HERNANDEZ',36,'LP1912'),
  @('11:15:53','12:34','11_ETCHEVERRY',79,'LP1912'),
  @('12:24:14','12:35','23_HERNANDEZ',11,'LP1912'),
  @('12:24:14','12:35','11_ETCHEVERRY',11,'LP1912'),
  @('10:50:37','12:36','27_EL RETIRO',106,'LP1912'),
  @('11:43:19','12:37','23_HERNANDEZ',54,'LP1912'),
  @('11:43:19','12:37','27_EL RETIRO',54,'LP1912'),
  @('10:50:37','12:38','17_179 Y 38',108,'LP1912'),
  @('11:15:53','12:40','10_OLMOS',85,'LP1912'),
  @('11:43:19','12:41','10_OLMOS',58,'LP1912'),
  @('12:44:05','12:44','17_ROMERO',0,'LP1912'),
  @('12:44:05','12:45','10_OLMOS',1,'LP1912'),
  @('12:44:05','12:45','16_SANTA ANA',1,'LP1912'),
  @('12:24:14','12:46','16_SANTA ANA',22,'LP1912'),
  @('11:15:53','12:46','17_ROMERO',91,'LP1912'),
  @('11:15:53','12:48','11_ETCHEVERRY',93,'LP1912'),
  @('12:24:14','12:49','11_ETCHEVERRY',25,'LP1912'),
  @('12:57:33','12:57','17_ROMERO',0,'LP1912'),
  @('12:57:33','12:57','16_SANTA ANA',0,'LP1912'),
  @('12:57:33','12:58','16_SANTA ANA',1,'LP1912'),
  @('11:15:53','13:02','15_ABASTO',107,'LP1912'),
  @('12:24:14','13:03','14_ABASTO',39,'LP1912'),
  @('11:15:53','13:06','16_P MOR-SANTA ANA',111,'LP1912'),
  @('12:57:33','13:06','23_HERNANDEZ',9,'LP1912'),
  @('12:24:14','13:07','16_P MOR-SANTA ANA',43,'LP1912'),
  @('11:15:53','13:13','215D_EL PATO',118,'LP1912'),
  @('12:57:33','13:14','11_ETCHEVERRY',17,'LP1912'),
  @('11:43:19','13:14','215D_EL PATO',91,'LP1912'),
  @('11:43:19','13:14','17_ROMERO',91,'LP1912'),
  @('12:44:05','13:15','11_ETCHEVERRY',31,'LP1912'),
  @('11:43:19','13:19','10_OLMOS',96,'LP1912'),
  @('12:44:05','13:20','10_OLMOS',36,'LP1912'),
  @('11:43:19','13:21','26_HERNANDEZ',98,'LP1912'),
  @('13:24:27','13:24','17_ROMERO',0,'LP1912'),
  @('13:24:27','13:25','16_SANTA ANA',1,'LP1912'),
  @('11:43:19','13:26','14_ABASTO',103,'LP1912'),
  @('11:43:19','13:26','15_ABASTO',103,'LP1912'),
  @('12:24:14','13:27','14_ABASTO',63,'LP1912'),
  @('11:58:46','13:29','17_ROMERO',91,'LP1912'),
  @('12:44:05','13:32','10_OLMOS',48,'LP1912'),
  @('12:57:33','13:34','23_HERNANDEZ',37,'LP1912'),
  @('13:24:27','13:36','15_ABASTO',12,'LP1912'),
  @('11:43:19','13:37','11_ETCHEVERRY',114,'LP1912'),
  @('12:24:14','13:40','23_HERNANDEZ',76,'LP1912'),
  @('12:44:05','13:41','23_HERNANDEZ',57,'LP1912'),
  @('11:58:46','13:46','17_ROMERO',108,'LP1912'),
  @('12:24:14','13:47','17_ROMERO',83,'LP1912'),
  @('12:57:33','13:50','11_ETCHEVERRY',53,'LP1912'),
  @('11:58:46','13:50','215A_EL PATO',112,'LP1912'),
  @('12:24:14','13:51','215A_EL PATO',87,'LP1912'),
  @('13:51:56','13:52','16_SANTA ANA',1,'LP1912'),
  @('13:51:56','13:55','225_GOMEZ',4,'LP1912'),
  @('11:58:46','13:56','16_P MOR-167 Y 521',118,'LP1912'),
  @('11:58:46','13:56','225_GOMEZ',118,'LP1912'),
  @('12:24:14','13:57','16_P MOR-167 Y 521',93,'LP1912'),
  @('13:51:56','13:59','17_ROMERO',8,'LP1912'),
  @('13:51:56','14:03','17_ROMERO',12,'LP1912'),
  @('13:51:56','14:04','11_ETCHEVERRY',13,'LP1912'),
  @('12:44:05','14:04','17_ROMERO',80,'LP1912'),
  @('13:24:27','14:04','23_HERNANDEZ',40,'LP1912'),
  @('12:44:05','14:05','23_HERNANDEZ',81,'LP1912'),
  @('13:24:27','14:05','11_ETCHEVERRY',41,'LP1912'),
  @('13:51:56','14:06','16_SANTA ANA',15,'LP1912'),
  @('13:51:56','14:11','15_ABASTO',20,'LP1912'),
  @('12:57:33','14:16','27_EL RETIRO',79,'LP1912'),
  @('12:24:14','14:17','27_EL RETIRO',113,'LP1912'),
  @('14:17:27','14:18','16_SANTA ANA',1,'LP1912'),
  @('13:51:56','14:19','215C_EL PATO',28,'LP1912'),
  @('13:51:56','14:20','26_HERNANDEZ',29,'LP1912'),
  @('12:24:14','14:20','215C_EL PATO',116,'LP1912'),
  @('12:24:14','14:21','26_HERNANDEZ',117,'LP1912'),
  @('14:17:27','14:28','15_ABASTO',11,'LP1912'),
  @('14:17:27','14:29','16_P MOR-SANTA ANA',12,'LP1912'),
  @('14:17:27','14:30','16_SANTA ANA',13,'LP1912'),
  @('14:17:27','14:34','23_HERNANDEZ',17,'LP1912'),
  @('14:42:26','14:42','14_ABASTO',0,'LP1912'),
  @('14:42:26','14:43','16_SANTA ANA',1,'LP1912'),
  @('13:51:56','14:44','14_ABASTO',53,'LP1912'),
  @('13:51:56','14:44','10_OLMOS',53,'LP1912'),
  @('12:57:33','14:45','14_ABASTO',108,'LP1912'),
  @('12:57:33','14:56','16_P MOR-SANTA ANA',119,'LP1912'),
  @('13:51:56','14:57','215B_EL PATO',66,'LP1912'),
  @('13:24:27','14:58','215B_EL PATO',94,'LP1912'),
  @('14:17:27','15:00','10_OLMOS',43,'LP1912'),
  @('13:24:27','15:00','81_EL PELIGRO',96,'LP1912'),
  @('13:51:56','15:04','10_OLMOS',73,'LP1912'),
  @('14:56:04','15:05','23_HERNANDEZ',9,'LP1912'),
  @('13:24:27','15:05','10_OLMOS',101,'LP1912'),
  @('14:42:26','15:06','16_SANTA ANA',24,'LP1912'),
  @('14:42:26','15:06','23_HERNANDEZ',24,'LP1912'),
  @('14:56:04','15:07','16_SANTA ANA',11,'LP1912'),
  @('14:42:26','15:10','10_OLMOS',28,'LP1912'),
  @('14:17:27','15:10','17_ROMERO',53,'LP1912'),
  @('13:51:56','15:13','11_ETCHEVERRY',82,'LP1912'),
  @('14:56:04','15:14','11_ETCHEVERRY',18,'LP1912'),
  @('14:56:04','15:17','16_SANTA ANA',21,'LP1912'),
  @('13:24:27','15:20','15_ABASTO',116,'LP1912'),
  @('13:51:56','15:21','26_HERNANDEZ',90,'LP1912'),
  @('13:24:27','15:22','26_HERNANDEZ',118,'LP1912'),
  @('15:21:47','15:24','10_OLMOS',3,'LP1912'),
  @('13:51:56','15:31','84_COLONIA URQUIZA-ESC 49',100,'LP1912'),
  @('14:17:27','15:32','84_COLONIA URQUIZA-ESC 49',75,'LP1912'),
  @('13:51:56','15:34','23_HERNANDEZ',103,'LP1912'),
  @('14:56:04','15:35','23_HERNANDEZ',39,'LP1912'),
  @('14:17:27','15:38','23_HERNANDEZ',81,'LP1912'),
  @('13:51:56','15:38','215A_EL PATO',107,'LP1912'),
  @('14:56:04','15:39','215A_EL PATO',43,'LP1912'),
  @('15:21:47','15:40','23_HERNANDEZ',19,'LP1912'),
  @('14:42:26','15:45','14_ABASTO',63,'LP1912'),
  @('14:56:04','15:46','14_ABASTO',50,'LP1912'),
  @('13:51:56','15:46','16_P MOR-167 Y 521',115,'LP1912'),
  @('14:56:04','15:47','16_P MOR-167 Y 521',51,'LP1912'),
  @('13:51:56','15:49','14_ABASTO',118,'LP1912'),
  @('15:21:47','15:54','16_SANTA ANA',33,'LP1912'),
  @('15:21:47','15:55','11_ETCHEVERRY',34,'LP1912'),
  @('15:21:47','15:55','17_ROMERO',34,'LP1912'),
  @('14:17:27','15:56','17_ROMERO',99,'LP1912'),
  @('14:17:27','15:56','27_EL RETIRO',99,'LP1912'),
  @('14:42:26','15:56','11_ETCHEVERRY',74,'LP1912'),
  @('14:17:27','15:57','11_ETCHEVERRY',100,'LP1912'),
  @('15:53:26','15:58','27_EL RETIRO',5,'LP1912'),
  @('14:56:04','16:01','27_EL RETIRO',65,'LP1912'),
  @('15:21:47','16:01','10_OLMOS',40,'LP1912'),
  @('15:21:47','16:02','16_SANTA ANA',41,'LP1912'),
  @('14:56:04','16:02','10_OLMOS',66,'LP1912'),
  @('15:21:47','16:04','23_HERNANDEZ',43,'LP1912'),
  @('14:42:26','16:05','27_EL RETIRO',83,'LP1912'),
  @('15:53:26','16:05','23_HERNANDEZ',12,'LP1912'),
  @('14:42:26','16:06','23_HERNANDEZ',84,'LP1912'),
  @('15:21:47','16:08','14_ABASTO',47,'LP1912'),
  @('14:56:04','16:09','14_ABASTO',73,'LP1912'),
  @('15:53:26','16:11','16_SANTA ANA',18,'LP1912'),
  @('15:21:47','16:14','225_C ROCA-H SUR',53,'LP1912'),
  @('14:17:27','16:15','225_C ROCA-H SUR',118,'LP1912'),
  @('15:21:47','16:19','215C_EL PATO',58,'LP1912'),
  @('16:14:21','16:19','16_SANTA ANA',5,'LP1912'),
  @('14:42:26','16:20','215C_EL PATO',98,'LP1912'),
  @('14:42:26','16:21','26_HERNANDEZ',99,'LP1912'),
  @('15:53:26','16:24','84_COLONIA URQUIZA-ESC 49',31,'LP1912'),
  @('16:14:21','16:27','16_SANTA ANA',13,'LP1912'),
  @('15:53:26','16:29','10_OLMOS',36,'LP1912'),
  @('14:42:26','16:30','15_ABASTO',108,'LP1912'),
  @('15:53:26','16:34','23_HERNANDEZ',41,'LP1912'),
  @('15:21:47','16:36','11_ETCHEVERRY',75,'LP1912'),
  @('15:21:47','16:39','17_ROMERO',78,'LP1912'),
  @('16:14:21','16:40','17_ROMERO',26,'LP1912'),
  @('15:21:47','16:42','16_P MOR-SANTA ANA',81,'LP1912'),
  @('14:56:04','16:43','225_GOMEZ',107,'LP1912'),
  @('14:56:04','16:43','16_P MOR-SANTA ANA',107,'LP1912'),
  @('16:14:21','16:43','23_HERNANDEZ',29,'LP1912'),
  @('15:21:47','16:48','15_ABASTO',87,'LP1912'),
  @('15:53:26','16:50','14_ABASTO',57,'LP1912'),
  @('15:21:47','16:56','17_179 Y 38',95,'LP1912'),
  @('16:14:21','16:57','10_OLMOS',43,'LP1912'),
  @('15:21:47','17:04','215A_EL PATO',103,'LP1912'),
  @('16:14:21','17:04','11_ETCHEVERRY',50,'LP1912'),
  @('16:14:21','17:05','23_HERNANDEZ',51,'LP1912'),
  @('15:53:26','17:21','26_HERNANDEZ',88,'LP1912'),
  @('15:53:26','17:24','84_COLONIA URQUIZA-ESC 49',91,'LP1912'),
  @('15:53:26','17:28','14_ABASTO',95,'LP1912'),
  @('16:14:21','17:37','27_EL RETIRO',83,'LP1912'),
  @('15:53:26','17:38','27_EL RETIRO',105,'LP1912'),
  @('16:14:21','17:38','17_ROMERO',84,'LP1912'),
  @('15:53:26','17:39','17_ROMERO',106,'LP1912'),
  @('15:53:26','17:40','215B_EL PATO',107,'LP1912'),
  @('15:53:26','17:50','16_P MOR-167 Y 521',117,'LP1912'),
  @('16:14:21','17:52','81_EL PELIGRO',98,'LP1912')
)
for ($i = 0; $i -lt $data1.Count; $i++) {
  $r = $i + 6
  $row = $data1[$i]
  $ws1.Cells.Item($r, 1).Value = $row[0]
  $ws1.Cells.Item($r, 2).Value = $row[1]
  $ws1.Cells.Item($r, 3).Value = $row[2]
  $ws1.Cells.Item($r, 4).Value = $row[3]
  $ws1.Cells.Item($r, 5).Value = $row[4]
}

# --- Sheet "LP1912-215" ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $updated

# --- Sheet "6203-6173" ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = $updated
$ws3.Range("A3").Value = 'Total filas: 47'
$data3 = @(
  @('00:04:05','00:05','215A_LA PLATA',1,'L6173'),
  @('03:54:40','05:43','215A_LA PLATA',109,'L6173'),
  @('04:18:06','05:44','215A_LA PLATA',86,'L6173'),
  @('04:53:50','06:08','215A_LA PLATA',75,'L6173'),
  @('04:18:06','06:09','215A_LA PLATA',111,'L6173'),
  @('04:53:50','06:32','215C_LA PLATA',99,'L6203'),
  @('04:40:32','06:33','215C_LA PLATA',113,'L6203'),
  @('05:18:23','07:00','215B_LP-P MOR-1 Y 57',102,'L6173'),
  @('05:47:32','07:35','215A_LA PLATA',108,'L6173'),
  @('06:52:38','08:06','215C_LA PLATA',74,'L6203'),
  @('06:37:24','08:07','215C_LA PLATA',90,'L6203'),
  @('07:14:27','08:14','215C_LA PLATA',60,'L6203'),
  @('07:44:08','08:17','215C_LA PLATA',33,'L6203'),
  @('07:57:27','08:18','215C_LA PLATA',21,'L6203'),
  @('06:37:24','08:30','215A_LA PLATA',113,'L6173'),
  @('06:52:38','08:34','215A_LA PLATA',102,'L6173'),
  @('07:14:27','08:35','215A_LA PLATA',81,'L6173'),
  @('08:33:47','08:37','215A_LA PLATA',4,'L6173'),
  @('08:16:48','09:08','215D_LA PLATA',52,'L6203'),
  @('07:14:27','09:09','215D_LA PLATA',115,'L6203'),
  @('08:33:47','09:10','215D_LA PLATA',37,'L6203'),
  @('08:54:42','09:11','215D_LA PLATA',17,'L6203'),
  @('08:47:19','09:12','215D_LA PLATA',25,'L6203'),
  @('08:16:48','10:02','215B_LP-P MOR-40 Y 115',106,'L6173'),
  @('08:47:19','10:03','215B_LP-P MOR-40 Y 115',76,'L6173'),
  @('10:50:37','10:53','215A_LA PLATA',3,'L6173'),
  @('09:25:30','10:54','215A_LA PLATA',89,'L6173'),
  @('10:50:37','11:13','215C_LA PLATA',23,'L6203'),
  @('09:25:30','11:14','215C_LA PLATA',109,'L6203'),
  @('11:15:53','11:15','215C_LA PLATA',0,'L6203'),
  @('10:11:11','12:04','215A_LA PLATA',113,'L6173'),
  @('11:15:53','12:53','215C_LA PLATA',98,'L6203'),
  @('11:43:19','12:54','215C_LA PLATA',71,'L6203'),
  @('11:58:46','13:30','215B_LP-P MOR-1 Y 57',92,'L6173'),
  @('11:43:19','13:31','215B_LP-P MOR-1 Y 57',108,'L6173'),
  @('13:51:56','14:08','215A_LA PLATA',17,'L6173'),
  @('12:24:14','14:09','215A_LA PLATA',105,'L6173'),
  @('13:51:56','14:52','215D_LA PLATA',61,'L6203'),
  @('12:57:33','14:53','215D_LA PLATA',116,'L6203'),
  @('13:51:56','15:33','215A_LA PLATA',102,'L6173'),
  @('14:17:27','15:34','215A_LA PLATA',77,'L6173'),
  @('15:21:47','16:13','215C_LA PLATA',52,'L6203'),
  @('14:17:27','16:14','215C_LA PLATA',117,'L6203'),
  @('15:21:47','16:52','215B_LP-P MOR-40 Y 115',91,'L6173'),
  @('14:56:04','16:53','215B_LP-P MOR-40 Y 115',117,'L6173'),
  @('15:21:47','17:14','215A_LA PLATA',113,'L6173'),
  @('16:14:21','18:04','16_P MOR-SANTA ANA',110,'L6203')
)
for ($i = 0; $i -lt $data3.Count; $i++) {
  $r = $i + 6
  $row = $data3[$i]
  $ws3.Cells.Item($r, 1).Value = $row[0]
  $ws3.Cells.Item($r, 2).Value = $row[1]
  $ws3.Cells.Item($r, 3).Value = $row[2]
  $ws3.Cells.Item($r, 4).Value = $row[3]
  $ws3.Cells.Item($r, 5).Value = $row[4]
}

Write-Host "Done: sheet1 rows=" $data1.Count " sheet3 rows=" $data3.Count